$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

function Set-TextValue($range, $text) {
    # Force a numeric-looking string to be stored as text, then strip the
    # residual "Text" number-format style so the cell keeps its original
    # (unstyled) appearance.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2: A2 becomes text "1" (was numeric 1), C2 gets numeric value 88 (was empty)
Set-TextValue $ws.Range("A2") "1"
$ws.Range("C2").Value = 88

# Row 3: A3 becomes text "2" (was numeric 2), B3 "lade" -> "Lade", C3 gets numeric value 91 (was empty)
Set-TextValue $ws.Range("A3") "2"
$ws.Range("B3").Value = "Lade"
$ws.Range("C3").Value = 91

# Row 4: A4 becomes text "3" (was numeric 3); B4/C4 unchanged
Set-TextValue $ws.Range("A4") "3"

# Row 5: A5 becomes text "1" (was numeric 1); C5 changes from 92 to 88
Set-TextValue $ws.Range("A5") "1"
$ws.Range("C5").Value = 88

# Row 6: A6 becomes text "2" (was numeric 2); B6 "lade" -> "Lade"; C6 changes from 88 to 91
Set-TextValue $ws.Range("A6") "2"
$ws.Range("B6").Value = "Lade"
$ws.Range("C6").Value = 91

# Row 7: A7 becomes text "3" (was numeric 3); B7/C7 unchanged (C7 stays 88)
Set-TextValue $ws.Range("A7") "3"

# Row 8: new row -> A8 text "5", B8 "New User", C8 numeric 77
Set-TextValue $ws.Range("A8") "5"
$ws.Range("B8").Value = "New User"
$ws.Range("C8").Value = 77
